# Port-level landings worksheet updates:
# - Row 37 becomes the "MONTEREY AREA TOTALS" / "Totals" summary row
#   (previously it duplicated the Santa Cruz port label).
# - Column A is widened to match column B (they are merged into a single
#   visual block), since the "MONTEREY AREA TOTALS" label no longer needs
#   its own narrower column.
# - The saved selection moves to the whole of column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 37 labels -------------------------------------------------
$ws.Range("A37").Value = "MONTEREY AREA TOTALS"
$ws.Range("B37").Value = "Totals"

# --- Match column A's width to column B's ---------------------------------
$ws.Columns(1).ColumnWidth = $ws.Columns(2).ColumnWidth

# --- Update the saved selection to span all of column A -------------------
$ws.Range("A1:A1048576").Select() | Out-Null
